$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73; existing rows 73-148 shift down to 74-149
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with the new data record
$ws.Cells.Item(73, 1).Value = 11
$ws.Cells.Item(73, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value = "Bíobío"
$ws.Cells.Item(73, 4).Value = 44629
$ws.Cells.Item(73, 5).Value = 8
$ws.Cells.Item(73, 6).Value = 100112003
$ws.Cells.Item(73, 7).Value = "Ajo"
$ws.Cells.Item(73, 8).Value = "Chino"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 140
$ws.Cells.Item(73, 11).Value = 19000
$ws.Cells.Item(73, 12).Value = 20000
$ws.Cells.Item(73, 13).Value = 19429
$ws.Cells.Item(73, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(73, 15).Value = "China"
$ws.Cells.Item(73, 16).Value = 1943
$ws.Cells.Item(73, 17).Value = 10
$ws.Cells.Item(73, 18).Value = "Hortaliza"
